$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 354.88235
$ws.Range("I33").Value = 354.88235
$ws.Range("K33").Value = 354.88235
$ws.Range("M33").Value = -125.88235

$ws.Range("H113").Value = 2780.3333
$ws.Range("I113").Value = 2261.2
$ws.Range("J113").Value = 2980
$ws.Range("K113").Value = 2261.2
$ws.Range("L113").Value = 2980
$ws.Range("M113").Value = 992.8000000000002
$ws.Range("N113").Value = -9488

$ws.Range("H121").Value = 1939
$ws.Range("I121").Value = 897.5
$ws.Range("J121").Value = 2633.3333
$ws.Range("K121").Value = 2692.5
$ws.Range("L121").Value = 7899.999899999999
$ws.Range("M121").Value = -945.5
$ws.Range("N121").Value = -11393.9999

$ws.Range("H137").Value = 1284.8649
$ws.Range("I137").Value = 1230.5883
$ws.Range("J137").Value = 1900
$ws.Range("K137").Value = 3691.7649
$ws.Range("L137").Value = 5700
$ws.Range("M137").Value = -1141.7649
$ws.Range("N137").Value = -10800

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1987.9231
$ws.Range("I61").Value = 1734.3
$ws.Range("J61").Value = 2833.3333
$ws.Range("K61").Value = 1734.3
$ws.Range("L61").Value = 2833.3333
$ws.Range("M61").Value = -1522.3
$ws.Range("N61").Value = -3257.3333

$ws.Range("H74").Value = 1222.6666
$ws.Range("I74").Value = 871.1111
$ws.Range("J74").Value = 1750
$ws.Range("K74").Value = 871.1111
$ws.Range("L74").Value = 1750
$ws.Range("M74").Value = 2.888900000000035
$ws.Range("N74").Value = -3498

$ws.Range("H77").Value = 1222.6666
$ws.Range("I77").Value = 871.1111
$ws.Range("J77").Value = 1750
$ws.Range("K77").Value = 4355.555499999999
$ws.Range("L77").Value = 8750
$ws.Range("M77").Value = 12.44450000000052
$ws.Range("N77").Value = -17486

$ws.Range("H132").Value = 4405.0444
$ws.Range("I132").Value = 5155.207
$ws.Range("J132").Value = 3045.375
$ws.Range("K132").Value = 15465.621
$ws.Range("L132").Value = 9136.125
$ws.Range("M132").Value = -12935.621
$ws.Range("N132").Value = -14196.125

$ws.Range("H136").Value = 1987.9231
$ws.Range("I136").Value = 1734.3
$ws.Range("J136").Value = 2833.3333
$ws.Range("K136").Value = 5202.9
$ws.Range("L136").Value = 8499.999899999999
$ws.Range("M136").Value = -2652.9
$ws.Range("N136").Value = -13599.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H88").Value = 38000
$ws.Range("J88").Value = 38000
$ws.Range("L88").Value = 38000
$ws.Range("N88").Value = -38812

$ws.Range("H91").Value = 38000
$ws.Range("J91").Value = 38000
$ws.Range("L91").Value = 38000
$ws.Range("N91").Value = -40808

$ws.Range("H134").Value = 2411.9211
$ws.Range("I134").Value = 2016.3214
$ws.Range("J134").Value = 3519.6
$ws.Range("K134").Value = 6048.9642
$ws.Range("L134").Value = 10558.8
$ws.Range("M134").Value = -3513.9642
$ws.Range("N134").Value = -15628.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2047.375
$ws.Range("I31").Value = 1672.5238
$ws.Range("J31").Value = 4671.3335
$ws.Range("K31").Value = 1672.5238
$ws.Range("L31").Value = 4671.3335
$ws.Range("M31").Value = -1377.5238
$ws.Range("N31").Value = -5261.3335

$ws.Range("H34").Value = 2047.375
$ws.Range("I34").Value = 1672.5238
$ws.Range("J34").Value = 4671.3335
$ws.Range("K34").Value = 1672.5238
$ws.Range("L34").Value = 4671.3335
$ws.Range("M34").Value = -1470.5238
$ws.Range("N34").Value = -5075.3335

$ws.Range("H58").Value = 674959
$ws.Range("I58").Value = 904657.1
$ws.Range("J58").Value = 2271.7144
$ws.Range("K58").Value = 904657.1
$ws.Range("L58").Value = 2271.7144
$ws.Range("M58").Value = -904454.1
$ws.Range("N58").Value = -2677.7144

$ws.Range("H94").Value = 1939.125
$ws.Range("J94").Value = 2185.5
$ws.Range("L94").Value = 2185.5
$ws.Range("N94").Value = -3087.5

$ws.Range("H132").Value = 411333.5
$ws.Range("I132").Value = 467351.97
$ws.Range("J132").Value = 5199.5
$ws.Range("K132").Value = 1402055.91
$ws.Range("L132").Value = 15598.5
$ws.Range("M132").Value = -1399525.91
$ws.Range("N132").Value = -20658.5

$ws.Range("H134").Value = 1747.6227
$ws.Range("I134").Value = 1208.6666
$ws.Range("J134").Value = 3249
$ws.Range("K134").Value = 3625.9998
$ws.Range("L134").Value = 9747
$ws.Range("M134").Value = -1090.9998
$ws.Range("N134").Value = -14817

$ws.Range("H136").Value = 674959
$ws.Range("I136").Value = 904657.1
$ws.Range("J136").Value = 2271.7144
$ws.Range("K136").Value = 2713971.3
$ws.Range("L136").Value = 6815.1432
$ws.Range("M136").Value = -2711421.3
$ws.Range("N136").Value = -11915.1432

$ws.Range("H139").Value = 43193.332
$ws.Range("J139").Value = 43193.332
$ws.Range("L139").Value = 43193.332
$ws.Range("N139").Value = -53473.332

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1791.58
$ws.Range("I132").Value = 1207.2368
$ws.Range("J132").Value = 3642
$ws.Range("K132").Value = 3621.7104
$ws.Range("L132").Value = 10926
$ws.Range("M132").Value = -1091.7104
$ws.Range("N132").Value = -15986

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 22798
$ws.Range("I18").Value = 2995
$ws.Range("J18").Value = 36000
$ws.Range("K18").Value = 2995
$ws.Range("L18").Value = 36000
$ws.Range("M18").Value = -2823
$ws.Range("N18").Value = -36344

$ws.Range("H20").Value = 15027752
$ws.Range("J20").Value = 11006
$ws.Range("L20").Value = 11006
$ws.Range("N20").Value = -11458

$ws.Range("H22").Value = 709.1111
$ws.Range("I22").Value = 390
$ws.Range("J22").Value = 749
$ws.Range("K22").Value = 390
$ws.Range("L22").Value = 749
$ws.Range("M22").Value = -95
$ws.Range("N22").Value = -1339

$ws.Range("H27").Value = 709.1111
$ws.Range("I27").Value = 390
$ws.Range("J27").Value = 749
$ws.Range("K27").Value = 390
$ws.Range("L27").Value = 749
$ws.Range("M27").Value = -283
$ws.Range("N27").Value = -963

$ws.Range("H46").Value = 693.5484
$ws.Range("I46").Value = 571.4286
$ws.Range("J46").Value = 1833.3334
$ws.Range("K46").Value = 571.4286
$ws.Range("L46").Value = 1833.3334
$ws.Range("M46").Value = -383.4286
$ws.Range("N46").Value = -2209.3334

$ws.Range("H132").Value = 7224.091
$ws.Range("I132").Value = 8983.333000000001
$ws.Range("J132").Value = 5113
$ws.Range("K132").Value = 26949.999
$ws.Range("L132").Value = 15339
$ws.Range("M132").Value = -24419.999
$ws.Range("N132").Value = -20399

$ws.Range("H136").Value = 24636050
$ws.Range("I136").Value = 29413060
$ws.Range("J136").Value = 1433429.2
$ws.Range("K136").Value = 88239180
$ws.Range("L136").Value = 4300287.6
$ws.Range("M136").Value = -88236630
$ws.Range("N136").Value = -4305387.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2417.5715
$ws.Range("I132").Value = 1785.3077
$ws.Range("J132").Value = 2965.5334
$ws.Range("K132").Value = 5355.9231
$ws.Range("L132").Value = 8896.600199999999
$ws.Range("M132").Value = -2825.9231
$ws.Range("N132").Value = -13956.6002

$ws.Range("H136").Value = 1827.24
$ws.Range("I136").Value = 1633.25
$ws.Range("J136").Value = 2172.111
$ws.Range("K136").Value = 4899.75
$ws.Range("L136").Value = 6516.333
$ws.Range("M136").Value = -2349.75
$ws.Range("N136").Value = -11616.333
